$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new description for the "network_keep_alive" command (row 14) ---
$ws.Range("D14").Value = "Keep alive message from CC2650 to gateway"

# --- Insert a new blank row at 36 (pushes "Gateway to sink UART" and the
#     code line below it down by one row) ---
$ws.Rows.Item(36).Insert()

# --- Fill in the previously-blank note cell (now A35) with the new note text ---
$ws.Range("A35").Value = "Note:  The n in (n) behind elements of the payload indicates the size in bytes of that element"

# Style: bold, size 12 like the heading "Note: " style used above it (A34),
# but use Text number format (same as before) for the cell.
$ws.Range("A35").NumberFormat = "@"
$ws.Range("A35").Font.Name = "Liberation Sans"
$ws.Range("A35").Font.Size = 12
$ws.Range("A35").Font.Bold = $true
$ws.Range("A35").Font.Color = 0

# The second run ("The n in (n) behind ... that element") is regular weight.
$run2 = $ws.Range("A35").Characters(8, 200)
$run2.Font.Bold = $false
$run2.Font.Name = "Liberation Sans"
$run2.Font.Size = 12
$run2.Font.Color = 0

# --- Restore the selection to C11 (matches the saved cursor position) ---
$ws.Range("C11").Select()
